$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.330.21"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.073.99"
$ws.Range("E3").Value = "  +4.51%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.16%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.34%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("D13").Value = "2.379.68"
$ws.Range("E13").Value = "  +4.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.97%  "

$ws.Range("D18").Value = "2.070.74"
$ws.Range("E18").Value = "  +4.09%  "

$ws.Range("D19").Value = "37.313.06"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +19.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").Value = "0.0₃0812"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.97%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("E25").Value = "  +3.07%  "

$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("E29").Value = "  +6.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("E31").Value = "  +5.95%  "

$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0622"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("E35").Value = "  +8.49%  "

$ws.Range("E36").Value = "  +3.60%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "

$ws.Range("E41").Value = "  -2.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +25.19%  "

$ws.Range("E43").Value = "  +7.77%  "

$ws.Range("D44").Value = "1.472.23"
$ws.Range("E44").Value = "  +3.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.16%  "

$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.54%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.90%  "
